$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: strip leftover (unused) cell styles from A5/C5 ---
$ws.Cells.Item(5,1).Style = "Normal"
$ws.Cells.Item(5,3).Style = "Normal"

# --- Row 12: was a styled/blank-padded row (SCuM QFN / U1 / Custom). ---
# Drop the whole 12:14 block (clears stray styled empty cells and the
# stale J1/J2 26-pos + J3 10-pos connector rows) then rebuild it clean.
$ws.Range("A12:P14").Delete()

$ws.Cells.Item(12,1).Value = "SCuM QFN"
$ws.Cells.Item(12,3).Value = "U1"
$ws.Cells.Item(12,5).Value = "Custom"
$ws.Cells.Item(12,9).Value = 1

# --- Row 13: J1, 22-position header ---
$ws.Cells.Item(13,1).Value = "CONN HEADER VERT 22POS 2.54MM"
$ws.Cells.Item(13,2).Value = "Connector Header Through Hole 22 position 0.100"" (2.54mm)"
$ws.Cells.Item(13,3).Value = "J1"
$ws.Cells.Item(13,4).Value = "Sullins Connector Solutions"
$ws.Cells.Item(13,5).Value = "PRPC011DAAN-RC"
$ws.Cells.Item(13,6).Value = "Digi-Key"
$ws.Cells.Item(13,8).Value = "https://www.digikey.com/en/products/detail/sullins-connector-solutions/PRPC011DAAN-RC/2775283"
$ws.Cells.Item(13,9).Value = 1

# --- Row 14: J2, 20-position header ---
$ws.Cells.Item(14,1).Value = "CONN HEADER VERT 20POS 2.54MM"
$ws.Cells.Item(14,2).Value = "Connector Header Through Hole 20 position 0.100"" (2.54mm)"
$ws.Cells.Item(14,3).Value = "J2"
$ws.Cells.Item(14,4).Value = "Sullins Connector Solutions"
$ws.Cells.Item(14,5).Value = "PRPC010DAAN-RC"
$ws.Cells.Item(14,6).Value = "Digi-Key"
$ws.Cells.Item(14,8).Value = "https://www.digikey.com/en/products/detail/sullins-connector-solutions/PRPC010DAAN-RC/2775284"
$ws.Cells.Item(14,9).Value = 1

# --- Row 15 (new row): J3, 10-position header ---
$ws.Cells.Item(15,1).Value = "CONN HEADER VERT 10POS 2.54MM"
$ws.Cells.Item(15,2).Value = "Connector Header Through Hole 10 position 0.100"" (2.54mm)"
$ws.Cells.Item(15,3).Value = "J3"
$ws.Cells.Item(15,4).Value = "Sullins Connector Solutions"
$ws.Cells.Item(15,5).Value = "PRPC010SAAN-RC"
$ws.Cells.Item(15,6).Value = "Digi-Key"
$ws.Cells.Item(15,8).Value = "https://www.digikey.com/en/products/detail/sullins-connector-solutions/PRPC010SAAN-RC/2775244"
$ws.Cells.Item(15,9).Value = 1

# Hyperlinks get (re-)added in this order -- matches the rId ordering (H14,
# H13, H15) left behind in the saved workbook.
$ws.Hyperlinks.Add($ws.Cells.Item(14,8), "https://www.digikey.com/en/products/detail/sullins-connector-solutions/PRPC010DAAN-RC/2775284")
$ws.Cells.Item(14,8).Style = $ws.Cells.Item(11,8).Style

$ws.Hyperlinks.Add($ws.Cells.Item(13,8), "https://www.digikey.com/en/products/detail/sullins-connector-solutions/PRPC011DAAN-RC/2775283")
$ws.Cells.Item(13,8).Style = $ws.Cells.Item(11,8).Style

$ws.Hyperlinks.Add($ws.Cells.Item(15,8), "https://www.digikey.com/en/products/detail/sullins-connector-solutions/PRPC010SAAN-RC/2775244")
$ws.Cells.Item(15,8).Style = $ws.Cells.Item(11,8).Style

# --- Selection left where the author's cursor ended up ---
$null = $ws.Range("B23").Select()
